# Actualización desde MV -datos-
# Adds four new daily rows (13-09-2021 .. 16-09-2021) to the bottom of the
# "Swap promedio cámara 2021 - Diaria" sheet, mirroring the pattern of the
# existing rows (column A = date label, columns B:M = numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "13-09-2021"; Values = @(3.94, 4.23, 4.42, 4.59, 5.09, -1.33, -0.07000000000000001, 0.31, 0.75, 1.01, 1.67, 1.86) },
    @{ Date = "14-09-2021"; Values = @(3.9,  4.16, 4.37, 4.54, 5.03, -1.29, -0.08, 0.29, 0.72, 1,    1.6,  1.78) },
    @{ Date = "15-09-2021"; Values = @(3.87, 4.11, 4.3,  4.46, 4.95, -1.24, -0.08, 0.29, 0.71, 0.98, 1.55, 1.71) },
    @{ Date = "16-09-2021"; Values = @(3.92, 4.16, 4.34, 4.5,  4.98, -1.15, $null, 0.36, 0.78, 1.03, 1.59, 1.75) }
)

$startRow = 180

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $entry.Date

    for ($c = 0; $c -lt $entry.Values.Count; $c++) {
        $val = $entry.Values[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 2).Value = $val
        }
    }
}
